{"js": "// Replace \"ADMINISTRATOR\" with \"MANAGER\" in the \"CARETAKER / OFFICE ADMINISTRATOR \"\n// job-title line of the Experience section, keeping the same run formatting\n// (Nunito font, size 21) used by the rest of that line.\n\nconst body = context.document.body;\nconst results = body.search(\"CARETAKER / OFFICE ADMINISTRATOR \", { matchCase: true });\nresults.load(\"items\");\nawait context.sync();\n\nif (results.items.length === 0) {\n  throw new Error(\"Could not find target text 'CARETAKER / OFFICE ADMINISTRATOR '\");\n}\n\nconst target = results.items[0];\n\n// Replace the whole line's text with the updated title, matching the net\n// textual result of the original diff (which re-split the single run into\n// \"CARETAKER / OFFICE \" + \"MANAGER\" + \" \", all sharing identical formatting).\ntarget.insertText(\"CARETAKER / OFFICE MANAGER \", \"Replace\");\nawait context.sync();\n", "ps1": "$d = $word.ActiveDocument\n\n# Update the job title \"CARETAKER / OFFICE ADMINISTRATOR \" to\n# \"CARETAKER / OFFICE MANAGER \" in the Experience section, preserving the\n# existing run formatting (Nunito font, size 21) already on that text.\n$find = $d.Content.Find\n$find.ClearFormatting()\n$find.Text = \"CARETAKER / OFFICE ADMINISTRATOR \"\n$find.Replacement.ClearFormatting()\n$find.Replacement.Text = \"CARETAKER / OFFICE MANAGER \"\n$find.Forward = $true\n$find.Wrap = 1\n$find.Format = $false\n$find.MatchCase = $true\n$find.MatchWholeWord = $false\n$find.MatchWildcards = $false\n$find.Execute([ref]$find.Text, [ref]$find.MatchCase, [ref]$find.MatchWholeWord, [ref]$find.MatchWildcards, $null, $null, [ref]$find.Forward, [ref]$find.Wrap, $null, [ref]$find.Replacement.Text, 2)\n"}
